$d = $word.ActiveDocument

$old = "Re: Submission of Original Research Article >=" + [char]0x20AC + [char]0x201C
$new = "Re: Submission of Original Research Article " + [char]0x20AC + [char]0x201C
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 0" }

$old = ">=" + [char]0x20AC + [char]0x00A2 + " All models achieved exceptional reproducibility (99-100% consistency)"
$new = "" + [char]0x20AC + [char]0x00A2 + " All models achieved exceptional reproducibility (99-100% consistency)"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 1" }

$old = ">=" + [char]0x20AC + [char]0x00A2 + " Diagnostic accuracy remained at chance level (~50%)"
$new = "" + [char]0x20AC + [char]0x00A2 + " Diagnostic accuracy remained at chance level (~50%)"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 2" }

$old = ">=" + [char]0x20AC + [char]0x00A2 + " The consistency-accuracy gap reached ~50 percentage points"
$new = "" + [char]0x20AC + [char]0x00A2 + " The consistency-accuracy gap reached ~50 percentage points"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 3" }

$old = ">=" + [char]0x20AC + [char]0x00A2 + " Models showed systematic bias toward positive diagnosis (49-51 false positives vs 0-1 false negatives)"
$new = "" + [char]0x2022 + " Models showed systematic bias toward positive diagnosis (49-51 false positives vs 0-1 false negatives)"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 4" }

$old = ">=" + [char]0x20AC + [char]0x00A2 + " Prompt engineering had minimal impact (<3% prediction change)"
$new = "" + [char]0x20AC + [char]0x00A2 + " Prompt engineering had minimal impact (<3% prediction change)"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 5" }

$old = ">=" + [char]0x20AC + [char]0x00A2 + " Error patterns were highly systematic across all three models"
$new = "" + [char]0x20AC + [char]0x00A2 + " Error patterns were highly systematic across all three models"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 6" }

$old = ">=" + [char]0x20AC + [char]0x00A2 + " Aligns with the journal" + [char]0x2019 + "s focus on AI in medicine and clinical decision support"
$new = "" + [char]0x20AC + [char]0x00A2 + " Aligns with the journal" + [char]0x2019 + "s focus on AI in medicine and clinical decision support"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 7" }

$old = ">=" + [char]0x20AC + [char]0x00A2 + " Addresses timely concerns about LLM reliability in healthcare"
$new = "" + [char]0x20AC + [char]0x00A2 + " Addresses timely concerns about LLM reliability in healthcare"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 8" }

$old = ">=" + [char]0x20AC + [char]0x00A2 + " Provides rigorous empirical evidence with immediate clinical implications"
$new = "" + [char]0x20AC + [char]0x00A2 + " Provides rigorous empirical evidence with immediate clinical implications"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 9" }

$old = ">=" + [char]0x20AC + [char]0x00A2 + " Appeals to diverse readership (clinicians, AI researchers, policymakers)"
$new = "" + [char]0x20AC + [char]0x00A2 + " Appeals to diverse readership (clinicians, AI researchers, policymakers)"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 10" }

$old = ">=" + [char]0x20AC + [char]0x00A2 + " Contributes to ongoing dialogue about responsible AI in medicine"
$new = "" + [char]0x20AC + [char]0x00A2 + " Contributes to ongoing dialogue about responsible AI in medicine"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 11" }

$old = "=" + [char]0x20AC + [char]0x00A2 + " This manuscript represents original work not previously published or under consideration elsewhere"
$new = "" + [char]0x20AC + [char]0x00A2 + " This manuscript represents original work not previously published or under consideration elsewhere"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 12" }

$old = "=" + [char]0x20AC + [char]0x00A2 + " A preprint version is available on medRxiv for community feedback and rapid dissemination"
$new = "" + [char]0x20AC + [char]0x00A2 + " A preprint version is available on medRxiv for community feedback and rapid dissemination"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 13" }

$old = "=" + [char]0x20AC + [char]0x00A2 + " All authors have approved the manuscript and agree with submission to JAMIA"
$new = "" + [char]0x20AC + [char]0x00A2 + " All authors have approved the manuscript and agree with submission to JAMIA"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 14" }

$old = "=" + [char]0x20AC + [char]0x00A2 + " We have no conflicts of interest to declare"
$new = "" + [char]0x20AC + [char]0x00A2 + " We have no conflicts of interest to declare"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 15" }

$old = "=" + [char]0x20AC + [char]0x00A2 + " The study used publicly available de-identified data and did not require IRB approval"
$new = "" + [char]0x20AC + [char]0x00A2 + " The study used publicly available de-identified data and did not require IRB approval"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 16" }

$old = "=" + [char]0x20AC + [char]0x00A2 + " All data, code, and analysis scripts will be made publicly available upon acceptance"
$new = "" + [char]0x20AC + [char]0x00A2 + " All data, code, and analysis scripts will be made publicly available upon acceptance"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 17" }

$old = "Dwi Anggriani >=" + [char]0x20AC + [char]0x201C + " Institut Sains Teknologi dan Kesehatan " + [char]0x2019 + "Aisyiyah Kendari"
$new = "Dwi Anggriani " + [char]0x20AC + [char]0x201C + " Institut Sains Teknologi dan Kesehatan " + [char]0x2019 + "Aisyiyah Kendari"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 18" }

$old = "Muhammad Atnang >=" + [char]0x20AC + [char]0x201C + " Institut Sains Teknologi dan Kesehatan " + [char]0x2019 + "Aisyiyah Kendari"
$new = "Muhammad Atnang " + [char]0x20AC + [char]0x201C + " Institut Sains Teknologi dan Kesehatan " + [char]0x2019 + "Aisyiyah Kendari"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 19" }

$old = "Kartini Aprilia Pratiwi Nuzry >=" + [char]0x20AC + [char]0x201C + " Institut Sains Teknologi dan Kesehatan " + [char]0x2019 + "Aisyiyah Kendari"
$new = "Kartini Aprilia Pratiwi Nuzry " + [char]0x20AC + [char]0x201C + " Institut Sains Teknologi dan Kesehatan " + [char]0x2019 + "Aisyiyah Kendari"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 20" }

$old = "=" + [char]0x20AC + [char]0x00A2 + " First systematic evaluation of LLM consistency versus accuracy in medical diagnosis"
$new = "" + [char]0x20AC + [char]0x00A2 + " First systematic evaluation of LLM consistency versus accuracy in medical diagnosis"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 21" }

$old = "=" + [char]0x20AC + [char]0x00A2 + " 1,200 predictions from three state-of-the-art models with rigorous checkpoint system"
$new = "" + [char]0x20AC + [char]0x00A2 + " 1,200 predictions from three state-of-the-art models with rigorous checkpoint system"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 22" }

$old = "=" + [char]0x20AC + [char]0x00A2 + " 99-100% consistency but only 50% accuracy >=" + [char]0x20AC + [char]0x201C + " unprecedented 50-point gap"
$new = "" + [char]0x20AC + [char]0x00A2 + " 99-100% consistency but only 50% accuracy " + [char]0x20AC + [char]0x201C + " unprecedented 50-point gap"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 23" }

$old = "=" + [char]0x20AC + [char]0x00A2 + " Systematic positive diagnosis bias (49-51 false positives, 0-1 false negatives)"
$new = "" + [char]0x2022 + " Systematic positive diagnosis bias (49-51 false positives, 0-1 false negatives)"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 24" }

$old = "=" + [char]0x20AC + [char]0x00A2 + " Prompt engineering had minimal effect, suggesting deep-rooted model behavior"
$new = "" + [char]0x20AC + [char]0x00A2 + " Prompt engineering had minimal effect, suggesting deep-rooted model behavior"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 25" }

$old = "=" + [char]0x20AC + [char]0x00A2 + " Recommends LLMs as supplementary tools, not primary diagnostic systems"
$new = "" + [char]0x20AC + [char]0x00A2 + " Recommends LLMs as supplementary tools, not primary diagnostic systems"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Output "NOT FOUND: 26" }

# Fix paragraph styles: BlockText -> BodyText (x2) and FirstParagraph -> BodyText (x1)
# These are the paragraphs starting with the DECLARATIONS bullets, the COMPETING INTERESTS
# heading, and the MANUSCRIPT HIGHLIGHTS bullets.
$d.Paragraphs.Item(21).Range.set_Style("Body Text")
$d.Paragraphs.Item(22).Range.set_Style("Body Text")
$d.Paragraphs.Item(32).Range.set_Style("Body Text")
